$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: fill in the previously-empty activity/hours for 13/09/2018
$ws.Range("B20").Value = "Discussion, planification et organisation de groupe"
$ws.Range("C20").Value = 1.5

# Rows 33-35 keep the same wording but the shared text is tweaked
# ("et débugging" -> ", débugging"). Re-set them first so the updated
# string occupies the lowest free shared-string slot.
$ws.Range("B33").Value = "Suite refactoring, débugging et Javadoc"
$ws.Range("B34").Value = "Suite refactoring, débugging et Javadoc"
$ws.Range("B35").Value = "Suite refactoring, débugging et Javadoc"

# Row 36: new distinct wording + updated hours
$ws.Range("B36").Value = "Suite et fin refactoring, débugging et Javadoc"
$ws.Range("C36").Value = 5

# Row 37: new distinct wording + hours entered
$ws.Range("B37").Value = "Relecture du rapport et du manuel d'utilisateur"
$ws.Range("C37").Value = 4.5

# Update the view: scroll so row 19 is at the top and select C38
$excel.Goto($ws.Range("A19"), $true)
$ws.Range("C38").Select()
